$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "[AI SPARK 챌린지] 수도관 누수 탐지 분류 문제 베이스라인 코드! (오디오 파일 전처리)"
$ws.Range("E4").Value = "https://teddylee777.github.io/kaggle/water-pipe-leak-classification"

$ws.Range("D44").Value = "Nimble: Parallel GPU Task Scheduling for DL - NIPS 논문 리뷰"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/108"

$ws.Range("D50").Value = "KIAS-SNU Winter Camp"
$ws.Range("E50").Value = "http://incredible.egloos.com/7530062"
